$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume 1h) keep their text formatting
# so that numeric/percentage-looking strings are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '315.41'
$ws.Range("E2").Value = '2.21%'
$ws.Range("D3").Value = '39.21'
$ws.Range("E3").Value = '-1.34%'
$ws.Range("D4").Value = '5.156'
$ws.Range("E4").Value = '0.56%'
$ws.Range("D5").Value = '0.08161'
$ws.Range("E5").Value = '0.34%'
$ws.Range("D6").Value = '1.973'
$ws.Range("E6").Value = '1.41%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '8.351'
$ws.Range("E7").Value = '2.87%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9359'
$ws.Range("E8").Value = '0.78%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1308'
$ws.Range("E9").Value = '-7.88%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1985'
$ws.Range("E10").Value = '2.79%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.09008'
$ws.Range("E11").Value = '-0.77%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03526'
$ws.Range("E12").Value = '0.47%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09735'
$ws.Range("E13").Value = '-0.73%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001407'
$ws.Range("E14").Value = '0.92%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006272'
$ws.Range("E15").Value = '7.06%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.611'
$ws.Range("E16").Value = '-7.78%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.403'
$ws.Range("E17").Value = '4.36%'
$ws.Range("D18").Value = '3.265'
$ws.Range("E18").Value = '-3.51%'
$ws.Range("E19").Value = '0.45%'
$ws.Range("E20").Value = '-0.17%'
$ws.Range("D21").Value = '5.023'
$ws.Range("E21").Value = '6.37%'
$ws.Range("D22").Value = '0.2489'
$ws.Range("E22").Value = '2.68%'
$ws.Range("D23").Value = '0.04386'
$ws.Range("E23").Value = '0.11%'
$ws.Range("D24").Value = '0.001247'
$ws.Range("E24").Value = '1.29%'
$ws.Range("D25").Value = '0.004751'
$ws.Range("E25").Value = '8.35%'
$ws.Range("D26").Value = '0.0003896'
$ws.Range("E26").Value = '199.36%'
$ws.Range("E27").Value = '-7.64%'
$ws.Range("D39").Value = '0.02242'
$ws.Range("E39").Value = '8.60%'
$ws.Range("D40").Value = '0.05231'
$ws.Range("E40").Value = '2.53%'
$ws.Range("D41").Value = '0.007748'
$ws.Range("E41").Value = '4.24%'
$ws.Range("D42").Value = '0.01032'
$ws.Range("E42").Value = '4.52%'
$ws.Range("D43").Value = '0.1397'
$ws.Range("E43").Value = '2.28%'
$ws.Range("D44").Value = '0.002103'
$ws.Range("E44").Value = '-1.37%'
$ws.Range("D45").Value = '0.009130'
$ws.Range("E45").Value = '-4.56%'
$ws.Range("D46").Value = '0.00006826'
$ws.Range("E46").Value = '6.96%'
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").Value = '0.05%'
$ws.Range("D48").Value = '0.003010'
$ws.Range("E48").Value = '10.84%'
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").Value = '0.05%'
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").Value = '0.05%'

Write-Host "Applied cell updates"
